$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.673.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.622.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '193.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.638'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.46%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.619.02'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.182'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.668'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.33%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '57.73'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000306'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.200.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.616.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.92%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.542.59'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.35%  '

$ws.Range("E20").Value = '  +2.78%  '

$ws.Range("E21").Value = '  +3.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '486.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.69%  '

$ws.Range("E25").Value = '  +2.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.41'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.122'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.86'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '611.82'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '40.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0834'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.409'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.148'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.98'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +16.52%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.317.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +19.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0455'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +14.24%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.139'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.54%  '

